# recursive update to sbol2 assembly untested
#
# Shifts the SynBioHub part URLs shown in row 3 (columns E:I) of the
# "temp_module" assembly table one step along, inserting the B0032
# part URL and pushing the rest of the sequence down so that the
# J23100 URL ends up in the last column (I3). Column D3 (BsaI) and
# the A4:B3-style well-position labels in row 4 are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value2 = "https://charmme.synbiohub.org/user/Gonza10V/CIDARMoCloKit/B0032/1"
$ws.Range("F3").Value2 = "https://charmme.synbiohub.org/user/Gonza10V/CIDARMoCloKit/B0015/1"
$ws.Range("G3").Value2 = "https://charmme.synbiohub.org/user/Gonza10V/CIDARMoCloKit/ComponentDefinition_dvk_backbone_core/1"
$ws.Range("H3").Value2 = "https://charmme.synbiohub.org/user/Gonza10V/CIDARMoCloKit/E0040m_gfp/1"
$ws.Range("I3").Value2 = "https://charmme.synbiohub.org/user/Gonza10V/CIDARMoCloKit/J23100/1"
